$wb = $excel.ActiveWorkbook

# --- INTER_SWITCH_LINKS sheet: add new edge-router ISL rows (39-42) ---
$ws = $wb.Worksheets.Item("INTER_SWITCH_LINKS")

# Insert 4 new formatted rows by cloning the format of the last existing
# data row (38), which carries the fill style used throughout this table,
# so the new rows visually match the existing ones.
for ($i = 0; $i -lt 4; $i++) {
    $ws.Rows.Item(38).Copy()
    $ws.Rows.Item(39).Insert(-4121)
}

# Row 39: sw-edge-001 <-> sw-spine-001 (port 7)
$ws.Range("K39").Value = "x3003"
$ws.Range("J39").Value = "sw-edge-001"
$ws.Range("L39").Value = "u40"
$ws.Range("O39").Value = 1
$ws.Range("P39").Value = "sw-spine-001"
$ws.Range("Q39").Value = "x3000"
$ws.Range("R39").Value = "u40"
$ws.Range("T39").Value = 7

# Row 40: sw-edge-002 <-> sw-spine-001 (port 8)
$ws.Range("J40").Value = "sw-edge-002"
$ws.Range("K40").Value = "x3003"
$ws.Range("L40").Value = "u40"
$ws.Range("O40").Value = 1
$ws.Range("P40").Value = "sw-spine-001"
$ws.Range("Q40").Value = "x3000"
$ws.Range("R40").Value = "u40"
$ws.Range("T40").Value = 8

# Row 41: sw-edge-001 <-> sw-spine-002 (port 7)
$ws.Range("J41").Value = "sw-edge-001"
$ws.Range("K41").Value = "x3003"
$ws.Range("L41").Value = "u40"
$ws.Range("O41").Value = 2
$ws.Range("P41").Value = "sw-spine-002"
$ws.Range("Q41").Value = "x3001"
$ws.Range("R41").Value = "u40"
$ws.Range("T41").Value = 7

# Row 42: sw-edge-002 <-> sw-spine-002 (port 8)
$ws.Range("J42").Value = "sw-edge-002"
$ws.Range("K42").Value = "x3003"
$ws.Range("L42").Value = "u40"
$ws.Range("O42").Value = 2
$ws.Range("P42").Value = "sw-spine-002"
$ws.Range("Q42").Value = "x3001"
$ws.Range("R42").Value = "u40"
$ws.Range("T42").Value = 8

# --- View state: INTER_SWITCH_LINKS becomes the active/selected tab,
#     with the selection resting on K45 (just past the new data) ---
$ws.Activate()
$ws.Range("K45").Select()

# COMPUTE_NODES is no longer the active tab; its own selection (J23)
# stays as-is, it's just not the focused sheet anymore.
$cn = $wb.Worksheets.Item("COMPUTE_NODES")
$cn.Range("J23").Select()

$ws.Activate()
